$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "MSIDN" / "MSI" values for the row belonging to customer 828959809
$ws.Range("C11").Value = "3016875982"
$ws.Range("D11").Value = "732111198172291"

# Update the selected cell (UI state) to D16
$ws.Range("D16").Select()
